# Daily attendance processing - reorder the "Recorded By" (column G) entries
# so that system/admin-style recorders are listed before human recorders.
#
# The canonical ordering (highest priority first) observed in the processed
# reports is:
#   admin@admin.com, System, system, backup@backdoor.com, <everyone else>
#
# Each G-cell holds a comma-separated list of recorder names/emails; this
# script re-sorts each list (stable sort) according to that priority order,
# leaving single-value cells (and already-correctly-ordered ones) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Get-RecorderRank($name) {
    if ($name.CompareTo("admin@admin.com") -eq 0) { return 0 }
    if ($name.CompareTo("System") -eq 0) { return 1 }
    if ($name.CompareTo("system") -eq 0) { return 2 }
    if ($name.CompareTo("backup@backdoor.com") -eq 0) { return 3 }
    return 999
}

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)   # Column G = "Recorded By"
    $value = $cell.Text

    if ([string]::IsNullOrEmpty($value)) {
        continue
    }

    $parts = @($value -split ",\s*")
    if ($parts.Count -lt 2) {
        continue
    }

    $items = @()
    for ($i = 0; $i -lt $parts.Count; $i++) {
        $items += [PSCustomObject]@{ Rank = (Get-RecorderRank $parts[$i]); Name = $parts[$i] }
    }

    $sorted = $items | Sort-Object -Property Rank

    $newParts = @()
    foreach ($o in $sorted) {
        $newParts += $o.Name
    }
    $newValue = $newParts -join ", "

    if ($newValue -ne $value) {
        $cell.Value = $newValue
    }
}
